# Add a "percentage" column to both sheets in the workbook:
#   - "PI hours": insert a new column D ("percentage") between "hours" (C)
#     and "dept" (which shifts from D to E).
#   - "dept hours": append a new column D ("percentage") after "hours" (C).
# The percentage is each row's share of hours out of the column total,
# expressed as a 0-100 number (not a 0-1 fraction).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "PI hours"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("PI hours")

# Insert a brand-new column at D, pushing the existing "dept" column to E.
$ws1.Columns.Item(4).Insert()

$ws1.Range("D1").Value = "percentage"

$ws1.Range("D2").Value = 51.19047619047619
$ws1.Range("D3").Value = 34.52380952380953
$ws1.Range("D4").Value = 14.28571428571429

# ---------------------------------------------------------------------
# Sheet 2: "dept hours"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("dept hours")

# This sheet has no existing column D, so a brand-new cell there starts
# out with the default (unstyled) format. Copy the header formatting
# (bold, border, centered) from the neighboring "hours" header first.
$ws2.Range("C1").Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("D1").Value = "percentage"

$ws2.Range("D2").Value = 42.63959390862944
$ws2.Range("D3").Value = 21.82741116751269
$ws2.Range("D4").Value = 14.72081218274112
$ws2.Range("D5").Value = 14.72081218274112
$ws2.Range("D6").Value = 6.091370558375634
